# Auto-generated Excel COM-interop script to apply scheduled-runner price/profit updates
# to the Bahamut_Profits workbook (sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1894907.4
$ws.Range("I6").Value = 1894907.4
$ws.Range("K6").Value = 5684722.199999999
$ws.Range("M6").Value = -5684610.199999999
$ws.Range("H40").Value = 38463652
$ws.Range("J40").Value = 43480390
$ws.Range("L40").Value = 43480390
$ws.Range("N40").Value = -43480740
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H106").Value = 2862.125
$ws.Range("I106").Value = 2670.0588
$ws.Range("J106").Value = 3328.5715
$ws.Range("K106").Value = 2670.0588
$ws.Range("L106").Value = 3328.5715
$ws.Range("M106").Value = -2039.0588
$ws.Range("N106").Value = -4590.5715
$ws.Range("H120").Value = 25000
$ws.Range("J120").Value = 25000
$ws.Range("L120").Value = 25000
$ws.Range("N120").Value = -34676
$ws.Range("H132").Value = 1431.9744
$ws.Range("I132").Value = 1431.9744
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4295.9232
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1765.9232
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 910.8333
$ws.Range("I137").Value = 863
$ws.Range("J137").Value = 1150
$ws.Range("K137").Value = 2589
$ws.Range("L137").Value = 3450
$ws.Range("M137").Value = -39
$ws.Range("N137").Value = -8550
$ws.Range("H138").Value = 4318.9897
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4318.9897
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12956.9691
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23236.9691

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1255.2
$ws.Range("I2").Value = 1408.6666
$ws.Range("K2").Value = 1408.6666
$ws.Range("M2").Value = -1295.6666
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H32").Value = 4654642
$ws.Range("I32").Value = 5872767
$ws.Range("J32").Value = 25767.6
$ws.Range("K32").Value = 5872767
$ws.Range("L32").Value = 25767.6
$ws.Range("M32").Value = -5872480
$ws.Range("N32").Value = -26341.6
$ws.Range("H116").Value = 1255.2
$ws.Range("I116").Value = 1408.6666
$ws.Range("K116").Value = 1408.6666
$ws.Range("M116").Value = 885.3334
$ws.Range("H132").Value = 1169.2373
$ws.Range("I132").Value = 936.2766
$ws.Range("J132").Value = 2081.6667
$ws.Range("K132").Value = 2808.8298
$ws.Range("L132").Value = 6245.000100000001
$ws.Range("M132").Value = -278.8298
$ws.Range("N132").Value = -11305.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1255.2
$ws.Range("I3").Value = 1408.6666
$ws.Range("K3").Value = 1408.6666
$ws.Range("M3").Value = -1294.6666
$ws.Range("H69").Value = 48000
$ws.Range("J69").Value = 48000
$ws.Range("L69").Value = 48000
$ws.Range("N69").Value = -49622
$ws.Range("H72").Value = 48000
$ws.Range("J72").Value = 48000
$ws.Range("L72").Value = 144000
$ws.Range("N72").Value = -152112
$ws.Range("H105").Value = 4708.36
$ws.Range("I105").Value = 4776.846
$ws.Range("J105").Value = 4634.1665
$ws.Range("K105").Value = 4776.846
$ws.Range("L105").Value = 4634.1665
$ws.Range("M105").Value = -3029.846
$ws.Range("N105").Value = -8128.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6000490
$ws.Range("I12").Value = 12000000
$ws.Range("J12").Value = 980
$ws.Range("K12").Value = 12000000
$ws.Range("L12").Value = 980
$ws.Range("M12").Value = -11999830
$ws.Range("N12").Value = -1320
$ws.Range("H68").Value = 16165.667
$ws.Range("J68").Value = 16165.667
$ws.Range("L68").Value = 16165.667
$ws.Range("N68").Value = -17663.667
$ws.Range("H71").Value = 16165.667
$ws.Range("J71").Value = 16165.667
$ws.Range("L71").Value = 48497.001
$ws.Range("N71").Value = -55985.001
$ws.Range("H132").Value = 1074.6034
$ws.Range("I132").Value = 850.5185
$ws.Range("J132").Value = 4099.75
$ws.Range("K132").Value = 2551.5555
$ws.Range("L132").Value = 12299.25
$ws.Range("M132").Value = -21.55549999999994
$ws.Range("N132").Value = -17359.25
$ws.Range("H134").Value = 1186.2245
$ws.Range("I134").Value = 969.32556
$ws.Range("J134").Value = 2740.6667
$ws.Range("K134").Value = 2907.97668
$ws.Range("L134").Value = 8222.000100000001
$ws.Range("M134").Value = -372.9766799999998
$ws.Range("N134").Value = -13292.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2193.4
$ws.Range("I4").Value = 45
$ws.Range("K4").Value = 135
$ws.Range("M4").Value = -23
$ws.Range("H5").Value = 1181.1818
$ws.Range("I5").Value = 899
$ws.Range("J5").Value = 1675
$ws.Range("K5").Value = 2697
$ws.Range("L5").Value = 5025
$ws.Range("M5").Value = -2585
$ws.Range("N5").Value = -5249
$ws.Range("H92").Value = 922.44446
$ws.Range("I92").Value = 283.66666
$ws.Range("J92").Value = 2200
$ws.Range("K92").Value = 850.9999799999999
$ws.Range("L92").Value = 6600
$ws.Range("M92").Value = 397.0000200000001
$ws.Range("N92").Value = -9096
$ws.Range("H122").Value = 1112113.4
$ws.Range("J122").Value = 2001456
$ws.Range("L122").Value = 18013104
$ws.Range("N122").Value = -18018004
$ws.Range("H131").Value = 852.58
$ws.Range("I131").Value = 408.75
$ws.Range("J131").Value = 891.1739
$ws.Range("K131").Value = 1226.25
$ws.Range("L131").Value = 2673.5217
$ws.Range("M131").Value = 3813.75
$ws.Range("N131").Value = -12753.5217
$ws.Range("H135").Value = 1181.1818
$ws.Range("I135").Value = 899
$ws.Range("J135").Value = 1675
$ws.Range("K135").Value = 8091
$ws.Range("L135").Value = 15075
$ws.Range("M135").Value = -5556
$ws.Range("N135").Value = -20145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 55000000
$ws.Range("I14").Value = 55000000
$ws.Range("K14").Value = 55000000
$ws.Range("M14").Value = -54999832
$ws.Range("H44").Value = 5742.2
$ws.Range("I44").Value = 111
$ws.Range("J44").Value = 7150
$ws.Range("K44").Value = 111
$ws.Range("L44").Value = 7150
$ws.Range("M44").Value = 485
$ws.Range("N44").Value = -8342
$ws.Range("H132").Value = 1685.3
$ws.Range("I132").Value = 1454.6285
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 4363.8855
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -1833.8855
$ws.Range("N132").Value = -14960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1483.6364
$ws.Range("I46").Value = 1552.5
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 1552.5
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -1364.5
$ws.Range("N46").Value = -1676
$ws.Range("H55").Value = 233.07408
$ws.Range("I55").Value = 164.17647
$ws.Range("J55").Value = 350.2
$ws.Range("K55").Value = 164.17647
$ws.Range("L55").Value = 350.2
$ws.Range("M55").Value = 8.823530000000005
$ws.Range("N55").Value = -696.2
$ws.Range("H122").Value = 2700
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5650
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 639.2632
$ws.Range("I107").Value = 920.5454999999999
$ws.Range("J107").Value = 252.5
$ws.Range("K107").Value = 2761.6365
$ws.Range("L107").Value = 757.5
$ws.Range("M107").Value = -841.6364999999996
$ws.Range("N107").Value = -4597.5
$ws.Range("H122").Value = 833.1667
$ws.Range("I122").Value = 851
$ws.Range("J122").Value = 797.5
$ws.Range("K122").Value = 2553
$ws.Range("L122").Value = 2392.5
$ws.Range("M122").Value = -103
$ws.Range("N122").Value = -7292.5
$ws.Range("H123").Value = 38909
$ws.Range("J123").Value = 38909
$ws.Range("L123").Value = 38909
$ws.Range("N123").Value = -48709
$ws.Range("H132").Value = 1190.0625
$ws.Range("I132").Value = 689.67566
$ws.Range("J132").Value = 2873.182
$ws.Range("K132").Value = 2069.02698
$ws.Range("L132").Value = 8619.545999999998
$ws.Range("M132").Value = 460.9730199999999
$ws.Range("N132").Value = -13679.546
